$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI edge-table rows (Lgi3 -> Stx1a) across Sending/Target
# clusters FAPs & sCs, per Dr Hou's advice (re-run with corrected cluster
# membership: 3 groups -> FAPs, sCs, ECs).
$rows = @(
  @("FAPs","Lgi3","Stx1a","ECs",  3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,0.9171583333333334,2.751475,0.1826479122984665,0.1826479122984664,1.408545170741667,12.676906536675,0.1635639863591576,0.1635639863591575),
  @("FAPs","Lgi3","Stx1a","FAPs", 3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,2.685518,8.056554,0.5348087002134706,0.5348087002134706,4.124340664378001,37.119065979402,0.478929333742017,0.478929333742017),
  @("FAPs","Lgi3","Stx1a","sCs",  3,1,1.535771,4.607313,0.8955152254457543,0.8955152254457543,3,1,1.418779,4.256337,0.2825433874880628,0.2825433874880629,2.178919643609,19.610276792481,0.2530219053445797,0.2530219053445798),
  @("sCs","Lgi3","Stx1a","ECs",   1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,0.9171583333333334,2.751475,0.1826479122984665,0.1826479122984664,0.164342850275,1.479085652475,0.0190839259393089,0.01908392593930889),
  @("sCs","Lgi3","Stx1a","FAPs",  1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,2.685518,8.056554,0.5348087002134706,0.5348087002134706,0.481209913866,4.330889224793999,0.05587936647145361,0.05587936647145361),
  @("sCs","Lgi3","Stx1a","sCs",   1,0.3333333333333333,0.179187,0.537561,0.1044847745542456,0.1044847745542456,3,1,1.418779,4.256337,0.2825433874880628,0.2825433874880629,0.254226752673,2.288040774057,0.02952148214348311,0.02952148214348312)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = 2 + $i
    $vals = $rows[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
